$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.1
$ws.Range("G2").Value = 1.12
$ws.Range("H2").Value = 100
$ws.Range("I2").Value = 160
$ws.Range("J2").Value = 10
$ws.Range("K2").Value = 11.5
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 2.16
$ws.Range("O2").Value = 1.83
$ws.Range("P2").Value = 1.21
$ws.Range("Q2").Value = 5.3
$ws.Range("R2").Value = 1.02
$ws.Range("S2").Value = 30
$ws.Range("T2").Value = 5.6
$ws.Range("U2").Value = 1.18
$ws.Range("V2").Value = 1.01
$ws.Range("W2").Value = 9
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 1000
$ws.Range("Z2").Value = 1000
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 2.22
$ws.Range("AC2").Value = 12
$ws.Range("AD2").Value = 180
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 3.9
$ws.Range("AG2").Value = 21
$ws.Range("AH2").Value = 310
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 15
$ws.Range("AK2").Value = 100
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 80
$ws.Range("AO2").Value = 1000
$ws.Range("F3").Value = 1.09
$ws.Range("G3").Value = 1.1
$ws.Range("H3").Value = 30
$ws.Range("I3").Value = 100
$ws.Range("J3").Value = 11.5
$ws.Range("K3").Value = 18
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 10
$ws.Range("O3").Value = 1.06
$ws.Range("P3").Value = 2.82
$ws.Range("Q3").Value = 1.25
$ws.Range("R3").Value = 1.5
$ws.Range("S3").Value = 2.32
$ws.Range("T3").Value = 1.01
$ws.Range("U3").Value = 1.02
$ws.Range("V3").Value = 1.02
$ws.Range("W3").Value = 8.800000000000001
$ws.Range("AB3").Value = 240
$ws.Range("AC3").Value = 790
$ws.Range("AF3").Value = 800
$ws.Range("AJ3").Value = 15
$ws.Range("AN3").Value = 4
$ws.Range("F4").Value = 1.04
$ws.Range("G4").Value = 980
$ws.Range("H4").Value = 1.09
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 3.5
$ws.Range("K4").Value = 3.8
$ws.Range("L4").Value = 1.03
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 3.95
$ws.Range("O4").Value = 1.04
$ws.Range("P4").Value = 1.04
$ws.Range("Q4").Value = 1.05
$ws.Range("R4").Value = 1.05
$ws.Range("S4").Value = 1.02
$ws.Range("T4").Value = 1.63
$ws.Range("U4").Value = 2.08
$ws.Range("V4").Value = 1.02
$ws.Range("W4").Value = 1.02
$ws.Range("X4").Value = 990
$ws.Range("Y4").Value = 990
$ws.Range("Z4").Value = 980
$ws.Range("AA4").Value = 120
$ws.Range("AB4").Value = 990
$ws.Range("AC4").Value = 990
$ws.Range("AD4").Value = 990
$ws.Range("AE4").Value = 980
$ws.Range("AF4").Value = 980
$ws.Range("AG4").Value = 990
$ws.Range("AH4").Value = 990
$ws.Range("AI4").Value = 980
$ws.Range("AJ4").Value = 980
$ws.Range("AK4").Value = 980
$ws.Range("AL4").Value = 980
$ws.Range("AM4").Value = 140
$ws.Range("AN4").Value = 980
$ws.Range("AO4").Value = 980
$ws.Range("F5").Value = 7.8
$ws.Range("G5").Value = 8
$ws.Range("H5").Value = 1.5
$ws.Range("I5").Value = 1.52
$ws.Range("J5").Value = 4.7
$ws.Range("K5").Value = 4.8
$ws.Range("L5").Value = 1.38
$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 4.3
$ws.Range("O5").Value = 1.29
$ws.Range("P5").Value = 2.1
$ws.Range("Q5").Value = 1.87
$ws.Range("R5").Value = 1.42
$ws.Range("S5").Value = 3.25
$ws.Range("T5").Value = 2.02
$ws.Range("U5").Value = 1.94
$ws.Range("V5").Value = 2.92
$ws.Range("W5").Value = 1.14
$ws.Range("X5").Value = 18
$ws.Range("Y5").Value = 7.8
$ws.Range("Z5").Value = 8
$ws.Range("AA5").Value = 12.5
$ws.Range("AB5").Value = 23
$ws.Range("AC5").Value = 10
$ws.Range("AD5").Value = 9.6
$ws.Range("AE5").Value = 15.5
$ws.Range("AF5").Value = 60
$ws.Range("AG5").Value = 28
$ws.Range("AH5").Value = 24
$ws.Range("AI5").Value = 36
$ws.Range("AJ5").Value = 260
$ws.Range("AK5").Value = 120
$ws.Range("AL5").Value = 100
$ws.Range("AM5").Value = 150
$ws.Range("AN5").Value = 140
$ws.Range("AO5").Value = 8
$ws.Range("F6").Value = 1.93
$ws.Range("G6").Value = 2.22
$ws.Range("H6").Value = 3.55
$ws.Range("I6").Value = 4.5
$ws.Range("J6").Value = 3.45
$ws.Range("K6").Value = 4.2
$ws.Range("L6").Value = 1.35
$ws.Range("M6").Value = 1.07
$ws.Range("N6").Value = 3.7
$ws.Range("O6").Value = 1.29
$ws.Range("P6").Value = 1.98
$ws.Range("Q6").Value = 1.77
$ws.Range("R6").Value = 1.41
$ws.Range("S6").Value = 2.92
$ws.Range("T6").Value = 1.67
$ws.Range("U6").Value = 2.12
$ws.Range("V6").Value = 1.29
$ws.Range("W6").Value = 1.82
$ws.Range("AG6").Value = 990
$ws.Range("F7").Value = 1.69
$ws.Range("G7").Value = 1.75
$ws.Range("H7").Value = 5.2
$ws.Range("I7").Value = 6.2
$ws.Range("J7").Value = 3.8
$ws.Range("K7").Value = 4.2
$ws.Range("L7").Value = 1.4
$ws.Range("M7").Value = 1.07
$ws.Range("N7").Value = 3.6
$ws.Range("O7").Value = 1.3
$ws.Range("P7").Value = 1.94
$ws.Range("Q7").Value = 1.88
$ws.Range("R7").Value = 1.35
$ws.Range("S7").Value = 3.3
$ws.Range("T7").Value = 1.83
$ws.Range("U7").Value = 1.92
$ws.Range("V7").Value = 1.2
$ws.Range("W7").Value = 2.32
$ws.Range("X7").Value = 16
$ws.Range("Y7").Value = 19
$ws.Range("Z7").Value = 46
$ws.Range("AA7").Value = 160
$ws.Range("AB7").Value = 8.4
$ws.Range("AC7").Value = 9
$ws.Range("AD7").Value = 22
$ws.Range("AE7").Value = 80
$ws.Range("AF7").Value = 10.5
$ws.Range("AG7").Value = 10.5
$ws.Range("AH7").Value = 23
$ws.Range("AI7").Value = 85
$ws.Range("AJ7").Value = 18.5
$ws.Range("AK7").Value = 19.5
$ws.Range("AL7").Value = 38
$ws.Range("AM7").Value = 130
$ws.Range("AN7").Value = 13
$ws.Range("AO7").Value = 95
